$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly entry (Fecha = 2021-12-20, serial 44550) is being recorded for
# "Camote" at Vega Central Mapocho de Santiago. It is inserted right after
# the existing header/fixed rows (row 11), pushing all the historical rows
# down by two, which is why every row from the old 12 through 58 now lives
# two rows further down (14..60) and two brand-new "Primera"/"Segunda" rows
# land at 12/13.
$ws.Rows("12:13").Insert()

$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C12").Value = "Metropolitana"
$ws.Range("D12").Value = 44550
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = 100114002
$ws.Range("G12").Value = "Camote"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 1060
$ws.Range("K12").Value = 11000
$ws.Range("L12").Value = 12000
$ws.Range("M12").Value = 11500
$ws.Range("N12").Value = "$/malla 18 kilos"
$ws.Range("O12").Value = "Perú"
$ws.Range("P12").Value = 639
$ws.Range("Q12").Value = 18
$ws.Range("R12").Value = "Hortaliza"

$ws.Range("A13").Value = 9
$ws.Range("B13").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C13").Value = "Metropolitana"
$ws.Range("D13").Value = 44550
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 100114002
$ws.Range("G13").Value = "Camote"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 430
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = 10000
$ws.Range("N13").Value = "$/malla 18 kilos"
$ws.Range("O13").Value = "Perú"
$ws.Range("P13").Value = 556
$ws.Range("Q13").Value = 18
$ws.Range("R13").Value = "Hortaliza"
